$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "...si cler que ie desirois..." -> "...si cler que je desirois..."
#    The lone "i" is already its own run (no color in rPr); just fix the
#    single character in place so no neighbouring runs are disturbed.
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("si cler que ", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$iPos1 = $r1.End
$d.Range($iPos1, $iPos1 + 1).Text = "j"

# ---------------------------------------------------------------------
# 2) "...ie desirois ie lay encore recuit..." -> "...je desirois je lay..."
#    Same fix for the second lone "i" run.
# ---------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("e desirois ", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$iPos2 = $r2.End
$d.Range($iPos2, $iPos2 + 1).Text = "j"

# ---------------------------------------------------------------------
# 3) "grattibroisse" -> "gratte" + "broisse" (spelling fix: i -> e),
#    with the corrected letter split into its own run, matching the
#    lone-letter-run style used elsewhere in this document (no explicit
#    color, just rtl).
# ---------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("grattibroisse", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$wordStart = $r3.Start

# Borrow the formatting of an existing bare (colorless) single-letter run
# so the split-off letter gets the same minimal rPr instead of inheriting
# the surrounding word's explicit color.
$bareSrc = $d.Content
$bareSrc.Find.Execute("si cler que ", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$bareRun = $d.Range($bareSrc.End, $bareSrc.End + 1)

$midChar = $d.Range($wordStart + 5, $wordStart + 6)
$midChar.FormattedText = $bareRun.FormattedText
$d.Range($wordStart + 5, $wordStart + 6).Text = "e"

# ---------------------------------------------------------------------
# 4) " enflammer. esta" -> " enflammer esta" (drop the stray period).
# ---------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("enflammer. esta", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$dotPos = $r4.Start + 9
$d.Range($dotPos, $dotPos + 1).Text = ""
